$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1395348837209302
$ws.Range("C2").Value = 0.6790697674418604
$ws.Range("J2").Value = 0.004651162790697674
$ws.Range("P2").Value = 0.1209302325581395
$ws.Range("S2").Value = 0.05581395348837209
$ws.Range("B3").Value = 0.02597402597402598
$ws.Range("C3").Value = 0.05194805194805195
$ws.Range("J3").Value = 0.01298701298701299
$ws.Range("P3").Value = 0.7792207792207793
$ws.Range("S3").Value = 0.1298701298701299
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.75
$ws.Range("B6").Value = 0.08648648648648649
$ws.Range("D6").Value = 0.01621621621621622
$ws.Range("F6").Value = 0.07567567567567568
$ws.Range("J6").Value = 0.1513513513513514
$ws.Range("O6").Value = 0.01081081081081081
$ws.Range("Q6").Value = 0.1243243243243243
$ws.Range("R6").Value = 0.1297297297297297
$ws.Range("S6").Value = 0.4054054054054054
$ws.Range("B7").Value = 0.08749999999999999
$ws.Range("D7").Value = 0.0125
$ws.Range("E7").Value = 0.00625
$ws.Range("F7").Value = 0.05
$ws.Range("J7").Value = 0.14375
$ws.Range("O7").Value = 0.0125
$ws.Range("Q7").Value = 0.1375
$ws.Range("R7").Value = 0.10625
$ws.Range("S7").Value = 0.44375
$ws.Range("B8").Value = 0.04931506849315068
$ws.Range("D8").Value = 0.02191780821917808
$ws.Range("E8").Value = 0.00273972602739726
$ws.Range("F8").Value = 0.05753424657534247
$ws.Range("J8").Value = 0.1068493150684932
$ws.Range("O8").Value = 0.00821917808219178
$ws.Range("Q8").Value = 0.1863013698630137
$ws.Range("R8").Value = 0.1150684931506849
$ws.Range("S8").Value = 0.4520547945205479
$ws.Range("B9").Value = 0.09345794392523364
$ws.Range("D9").Value = 0.01869158878504673
$ws.Range("E9").Value = 0.004672897196261682
$ws.Range("F9").Value = 0.05607476635514019
$ws.Range("J9").Value = 0.1261682242990654
$ws.Range("O9").Value = 0.01401869158878505
$ws.Range("Q9").Value = 0.1822429906542056
$ws.Range("R9").Value = 0.1121495327102804
$ws.Range("S9").Value = 0.3925233644859813
$ws.Range("B10").Value = 0.09306260575296109
$ws.Range("D10").Value = 0.01269035532994924
$ws.Range("E10").Value = 0.001692047377326565
$ws.Range("F10").Value = 0.07360406091370558
$ws.Range("J10").Value = 0.1015228426395939
$ws.Range("O10").Value = 0.01607445008460237
$ws.Range("Q10").Value = 0.1852791878172589
$ws.Range("R10").Value = 0.1192893401015228
$ws.Range("S10").Value = 0.3967851099830795
$ws.Range("G11").Value = 0.1446280991735537
$ws.Range("J11").Value = 0.1033057851239669
$ws.Range("K11").Value = 0.2066115702479339
$ws.Range("L11").Value = 0.5330578512396694
$ws.Range("S11").Value = 0.01239669421487603
$ws.Range("G12").Value = 0.7647058823529411
$ws.Range("J12").Value = 0.1397058823529412
$ws.Range("K12").Value = 0.01470588235294118
$ws.Range("L12").Value = 0.05882352941176471
$ws.Range("S12").Value = 0.02205882352941177
$ws.Range("F15").Value = 0.01234567901234568
$ws.Range("I15").Value = 0.06584362139917696
$ws.Range("J15").Value = 0.4156378600823045
$ws.Range("K15").Value = 0.06584362139917696
$ws.Range("M15").Value = 0.00411522633744856
$ws.Range("O15").Value = 0.09053497942386832
$ws.Range("S15").Value = 0.2345679012345679
$ws.Range("H16").Value = 0.155688622754491
$ws.Range("I16").Value = 0.08383233532934131
$ws.Range("J16").Value = 0.4910179640718563
$ws.Range("K16").Value = 0.08982035928143713
$ws.Range("M16").Value = 0.02395209580838323
$ws.Range("O16").Value = 0.09580838323353294
$ws.Range("S16").Value = 0.05988023952095808
$ws.Range("F17").Value = 0.01084010840108401
$ws.Range("H17").Value = 0.1707317073170732
$ws.Range("I17").Value = 0.1111111111111111
$ws.Range("J17").Value = 0.4552845528455284
$ws.Range("K17").Value = 0.06504065040650407
$ws.Range("M17").Value = 0.008130081300813009
$ws.Range("O17").Value = 0.06233062330623306
$ws.Range("S17").Value = 0.1165311653116531
$ws.Range("F18").Value = 0.01204819277108434
$ws.Range("H18").Value = 0.1606425702811245
$ws.Range("I18").Value = 0.1044176706827309
$ws.Range("J18").Value = 0.4819277108433735
$ws.Range("K18").Value = 0.07228915662650602
$ws.Range("M18").Value = 0.02409638554216868
$ws.Range("O18").Value = 0.07630522088353414
$ws.Range("S18").Value = 0.06827309236947791
$ws.Range("F19").Value = 0.01323918799646955
$ws.Range("H19").Value = 0.1844660194174757
$ws.Range("I19").Value = 0.1032656663724625
$ws.Range("J19").Value = 0.3874669020300088
$ws.Range("K19").Value = 0.09885260370697264
$ws.Range("M19").Value = 0.01412180052956752
$ws.Range("N19").Value = 0.00176522506619594
$ws.Range("O19").Value = 0.09002647837599294
$ws.Range("S19").Value = 0.1067961165048544
